# "Ordering into modules write and read"
#
# The sheet lists people keyed by a "module" id. Originally the sheet had an
# extra lettered column (B/C/D/A) sitting between the name columns and the
# two module-id columns; that column is dropped, and the remaining rows are
# re-ordered so that all rows belonging to the same module (now column D)
# are grouped together (module 1386825528 block first, then module
# 7285813456 block) instead of being interleaved by write/read order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the lettered column (old column D); E/F shift left into D/E.
$ws.Columns.Item(4).Delete()

# Group rows by module id (now column D), stable/ascending so the
# 1386825528 block (originally the "write" rows 7-16) comes before the
# 7285813456 block (originally the "read" rows 1-6).
$rng = $ws.Range("A1:E16")
$rng.Sort($ws.Range("D1"), 1)

# The row that used to close out the first (write) module block
# (47708660 / Anette / Myra) is retired; a new row is appended to close out
# the second (read) module block instead.
$ws.Rows.Item(10).Delete()

$ws.Cells.Item(16, 1).Value = 32039648
$ws.Cells.Item(16, 2).Value = "Luciana"
$ws.Cells.Item(16, 3).Value = "Isabella"
$ws.Cells.Item(16, 4).Value = 7285813456
$ws.Cells.Item(16, 5).Value = 7116301002
